$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.841.57'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.888.62'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7522'
$ws.Range("E5").Value = '  -2.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.38'
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3129'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.33'
$ws.Range("E9").Value = '  -1.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07128'
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08499'
$ws.Range("E11").Value = '  +5.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7603'
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.900.16'
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.376'
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.35'
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.140'
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.942.70'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.71'
$ws.Range("E18").Value = '  -2.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.14'
$ws.Range("E19").Value = '  -1.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007846'
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.143.00'
$ws.Range("E22").Value = '  -2.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.011'
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1591'
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.375'
$ws.Range("E26").Value = '  -0.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.10'
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.71'
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.032'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.482'
$ws.Range("E30").Value = '  +3.44%  '
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.515'
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.166'
$ws.Range("E33").Value = '  +2.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05422'
$ws.Range("E34").Value = '  -2.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.242'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7529'
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.004'
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.709'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01950'
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4471'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.101.95'
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.108'
$ws.Range("E43").Value = '  +1.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.66'
$ws.Range("E44").Value = '  -2.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8599'
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.734'
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.56'
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.860'
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.041'
$ws.Range("E50").Value = '  +1.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.040.61'
$ws.Range("E51").Value = '  +0.07%  '
